$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 5.305917
$ws.Range("H2").Value = 15.917751
$ws.Range("I2").Value = 0.4336744870332215
$ws.Range("J2").Value = 0.4336744870332215
$ws.Range("M2").Value = 0.6946430000000001
$ws.Range("N2").Value = 2.083929
$ws.Range("O2").Value = 0.1140293552421611
$ws.Range("P2").Value = 0.1140293552421611
$ws.Range("Q2").Value = 3.685718102631001
$ws.Range("R2").Value = 33.171462923679
$ws.Range("S2").Value = 0.04945162214137322
$ws.Range("T2").Value = 0.04945162214137321

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 5.305917
$ws.Range("H3").Value = 15.917751
$ws.Range("I3").Value = 0.4336744870332215
$ws.Range("J3").Value = 0.4336744870332215
$ws.Range("M3").Value = 5.373609333333333
$ws.Range("N3").Value = 16.120828
$ws.Range("O3").Value = 0.8821066470161785
$ws.Range("P3").Value = 0.8821066470161785
$ws.Range("Q3").Value = 28.511925113092
$ws.Range("R3").Value = 256.607326017828
$ws.Range("S3").Value = 0.3825471476533362
$ws.Range("T3").Value = 0.3825471476533361

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 5.305917
$ws.Range("H4").Value = 15.917751
$ws.Range("I4").Value = 0.4336744870332215
$ws.Range("J4").Value = 0.4336744870332215
$ws.Range("M4").Value = 0.02353866666666667
$ws.Range("N4").Value = 0.070616
$ws.Range("O4").Value = 0.00386399774166032
$ws.Range("P4").Value = 0.00386399774166032
$ws.Range("Q4").Value = 0.124894211624
$ws.Range("R4").Value = 1.124047904616
$ws.Range("S4").Value = 0.001675717238512066
$ws.Range("T4").Value = 0.001675717238512066

# Row 5
$ws.Range("I5").Value = 0.2066699405724794
$ws.Range("J5").Value = 0.2066699405724794
$ws.Range("M5").Value = 0.6946430000000001
$ws.Range("N5").Value = 2.083929
$ws.Range("O5").Value = 0.1140293552421611
$ws.Range("P5").Value = 0.1140293552421611
$ws.Range("Q5").Value = 1.756449051104334
$ws.Range("R5").Value = 15.808041459939
$ws.Range("S5").Value = 0.02356644007141559
$ws.Range("T5").Value = 0.02356644007141559

# Row 6
$ws.Range("I6").Value = 0.2066699405724794
$ws.Range("J6").Value = 0.2066699405724794
$ws.Range("M6").Value = 5.373609333333333
$ws.Range("N6").Value = 16.120828
$ws.Range("O6").Value = 0.8821066470161785
$ws.Range("P6").Value = 0.8821066470161785
$ws.Range("Q6").Value = 13.58751331912756
$ws.Range("R6").Value = 122.287619872148
$ws.Range("S6").Value = 0.1823049283174227
$ws.Range("T6").Value = 0.1823049283174227

# Row 7
$ws.Range("I7").Value = 0.2066699405724794
$ws.Range("J7").Value = 0.2066699405724794
$ws.Range("M7").Value = 0.02353866666666667
$ws.Range("N7").Value = 0.070616
$ws.Range("O7").Value = 0.00386399774166032
$ws.Range("P7").Value = 0.00386399774166032
$ws.Range("Q7").Value = 0.05951901729511112
$ws.Range("R7").Value = 0.535671155656
$ws.Range("S7").Value = 0.0007985721836411331
$ws.Range("T7").Value = 0.0007985721836411331

# Row 8
$ws.Range("G8").Value = 4.400310999999999
$ws.Range("H8").Value = 13.200933
$ws.Range("I8").Value = 0.359655572394299
$ws.Range("J8").Value = 0.359655572394299
$ws.Range("M8").Value = 0.6946430000000001
$ws.Range("N8").Value = 2.083929
$ws.Range("O8").Value = 0.1140293552421611
$ws.Range("P8").Value = 0.1140293552421611
$ws.Range("Q8").Value = 3.056645233973
$ws.Range("R8").Value = 27.509807105757
$ws.Range("S8").Value = 0.04101129302937233
$ws.Range("T8").Value = 0.04101129302937233

# Row 9
$ws.Range("G9").Value = 4.400310999999999
$ws.Range("H9").Value = 13.200933
$ws.Range("I9").Value = 0.359655572394299
$ws.Range("J9").Value = 0.359655572394299
$ws.Range("M9").Value = 5.373609333333333
$ws.Range("N9").Value = 16.120828
$ws.Range("O9").Value = 0.8821066470161785
$ws.Range("P9").Value = 0.8821066470161785
$ws.Range("Q9").Value = 23.64555225916933
$ws.Range("R9").Value = 212.809970332524
$ws.Range("S9").Value = 0.3172545710454195
$ws.Range("T9").Value = 0.3172545710454195

# Row 10
$ws.Range("G10").Value = 4.400310999999999
$ws.Range("H10").Value = 13.200933
$ws.Range("I10").Value = 0.359655572394299
$ws.Range("J10").Value = 0.359655572394299
$ws.Range("M10").Value = 0.02353866666666667
$ws.Range("N10").Value = 0.070616
$ws.Range("O10").Value = 0.00386399774166032
$ws.Range("P10").Value = 0.00386399774166032
$ws.Range("Q10").Value = 0.1035774538586667
$ws.Range("R10").Value = 0.9321970847279999
$ws.Range("S10").Value = 0.001389708319507121
$ws.Range("T10").Value = 0.001389708319507121
